# Glossary workbook update: add a "Name" column to the Greek Symbols sheet
# and fill in missing "Name" values ("-") on the Math Operators sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Greek Symbols" sheet: insert a new column B ("Name") that spells
#    out the Greek letter name for each symbol in column A.
# ---------------------------------------------------------------------
$wsGreek = $wb.Worksheets.Item("Greek Symbols")

$wsGreek.Columns("B").Insert()
$wsGreek.Columns("B:L").Hidden = $false

$wsGreek.Range("B2").Value = "Name"

$greekNames = @{
    3  = "alpha"
    4  = "beta"
    5  = "gamma"
    6  = "gamma dot"
    7  = "delta"
    8  = "Delta"
    9  = "epsilon"
    10 = "Epsilon"
    11 = "Epsilon_lambda"
    12 = "zeta"
    13 = "eta"
    14 = "theta"
    15 = "theta dot"
    16 = "theta double dot"
    17 = "kappa"
    18 = "lambda"
    19 = "mu"
    20 = "nu"
    21 = "xi"
    22 = "rho"
    23 = "sigma"
    24 = "sigma"
    25 = "sigma"
    26 = "sigma"
    27 = "sigma"
    28 = "tau"
    29 = "tau"
    30 = "phi"
    31 = "phi"
    32 = "chi"
    33 = "omega"
}

foreach ($row in 3..33) {
    $wsGreek.Range("B$row").Value = $greekNames[$row]
}

# Row 25 carries an explicit style on column A; mirror it on the new cell.
$wsGreek.Range("B25").Style = $wsGreek.Range("A25").Style

# ---------------------------------------------------------------------
# 2. "Math Operators" sheet: fill the blank "Name" cells with "-".
# ---------------------------------------------------------------------
$wsMath = $wb.Worksheets.Item("Math Operators")

foreach ($row in 4..26) {
    $wsMath.Range("B$row").Value = "-"
}

# ---------------------------------------------------------------------
# 3. View / selection tweaks to match the saved workbook state.
# ---------------------------------------------------------------------
$wsGreek.Application.ActiveWindow.ScrollRow = 22
$wsGreek.Range("C26").Select()

$wsMath.Activate()
$wsMath.Range("Q21").Select()

$wsLatin = $wb.Worksheets.Item("Latin Symbols")
$wsLatin.Range("D14").Select()

Write-Host "done"
